$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 40.294117
$ws.Cells.Item(2, 10).Value = 192.5
$ws.Cells.Item(2, 12).Value = 192.5
$ws.Cells.Item(2, 14).Value = -418.5
$ws.Cells.Item(28, 8).Value = 1536.2106
$ws.Cells.Item(28, 9).Value = 1421.1
$ws.Cells.Item(28, 10).Value = 1664.1111
$ws.Cells.Item(28, 11).Value = 1421.1
$ws.Cells.Item(28, 12).Value = 1664.1111
$ws.Cells.Item(28, 13).Value = -936.0999999999999
$ws.Cells.Item(28, 14).Value = -2634.1111
$ws.Cells.Item(29, 8).Value = 1425
$ws.Cells.Item(29, 9).Value = 887.5
$ws.Cells.Item(29, 11).Value = 2662.5
$ws.Cells.Item(29, 13).Value = -2381.5
$ws.Cells.Item(40, 8).Value = 4137.76
$ws.Cells.Item(40, 9).Value = 2984.7856
$ws.Cells.Item(40, 11).Value = 2984.7856
$ws.Cells.Item(40, 13).Value = -2809.7856
$ws.Cells.Item(107, 8).Value = 31948544
$ws.Cells.Item(107, 9).Value = 17311292
$ws.Cells.Item(107, 11).Value = 17311292
$ws.Cells.Item(107, 13).Value = -17309372
$ws.Cells.Item(132, 8).Value = 1068.909
$ws.Cells.Item(132, 9).Value = 958.13116
$ws.Cells.Item(132, 11).Value = 2874.39348
$ws.Cells.Item(132, 13).Value = -344.3934800000002
$ws.Cells.Item(138, 8).Value = 4956.985
$ws.Cells.Item(138, 10).Value = 8212.105
$ws.Cells.Item(138, 12).Value = 24636.315
$ws.Cells.Item(138, 14).Value = -34916.315

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(29, 8).Value = 2994
$ws.Cells.Item(29, 9).Value = 933
$ws.Cells.Item(29, 10).Value = 7116
$ws.Cells.Item(29, 11).Value = 933
$ws.Cells.Item(29, 12).Value = 7116
$ws.Cells.Item(29, 13).Value = -625
$ws.Cells.Item(29, 14).Value = -7732
$ws.Cells.Item(82, 8).Value = 23453.285
$ws.Cells.Item(82, 10).Value = 23453.285
$ws.Cells.Item(82, 12).Value = 23453.285
$ws.Cells.Item(82, 14).Value = -24175.285
$ws.Cells.Item(85, 8).Value = 23453.285
$ws.Cells.Item(85, 10).Value = 23453.285
$ws.Cells.Item(85, 12).Value = 23453.285
$ws.Cells.Item(85, 14).Value = -25949.285
$ws.Cells.Item(97, 8).Value = 3788643.2
$ws.Cells.Item(97, 9).Value = 681.6667
$ws.Cells.Item(97, 10).Value = 11905704
$ws.Cells.Item(97, 11).Value = 681.6667
$ws.Cells.Item(97, 12).Value = 11905704
$ws.Cells.Item(97, 13).Value = -185.6667
$ws.Cells.Item(97, 14).Value = -11906696
$ws.Cells.Item(102, 8).Value = 1005.8
$ws.Cells.Item(102, 9).Value = 1005.8
$ws.Cells.Item(102, 11).Value = 1005.8
$ws.Cells.Item(102, 13).Value = 616.2
$ws.Cells.Item(110, 8).Value = 22223676
$ws.Cells.Item(110, 9).Value = 1365.3636
$ws.Cells.Item(110, 10).Value = 83335030
$ws.Cells.Item(110, 11).Value = 1365.3636
$ws.Cells.Item(110, 12).Value = 83335030
$ws.Cells.Item(110, 13).Value = 679.6364000000001
$ws.Cells.Item(110, 14).Value = -83339120
$ws.Cells.Item(122, 8).Value = 2526.1943
$ws.Cells.Item(122, 9).Value = 1734.3334
$ws.Cells.Item(122, 10).Value = 4901.778
$ws.Cells.Item(122, 11).Value = 5203.0002
$ws.Cells.Item(122, 12).Value = 14705.334
$ws.Cells.Item(122, 13).Value = -2753.0002
$ws.Cells.Item(122, 14).Value = -19605.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 43272316
$ws.Cells.Item(107, 9).Value = 59212120
$ws.Cells.Item(107, 11).Value = 59212120
$ws.Cells.Item(107, 13).Value = -59210200

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(21, 8).Value = 1323
$ws.Cells.Item(21, 10).Value = 1323
$ws.Cells.Item(21, 12).Value = 1323
$ws.Cells.Item(21, 14).Value = -1793
$ws.Cells.Item(22, 8).Value = 327.66666
$ws.Cells.Item(22, 9).Value = 313.4
$ws.Cells.Item(22, 10).Value = 399
$ws.Cells.Item(22, 11).Value = 313.4
$ws.Cells.Item(22, 12).Value = 399
$ws.Cells.Item(22, 13).Value = 36.60000000000002
$ws.Cells.Item(22, 14).Value = -1099
$ws.Cells.Item(26, 8).Value = 14000
$ws.Cells.Item(26, 9).Value = 1000
$ws.Cells.Item(26, 10).Value = 18333.334
$ws.Cells.Item(26, 11).Value = 1000
$ws.Cells.Item(26, 12).Value = 18333.334
$ws.Cells.Item(26, 13).Value = -713
$ws.Cells.Item(26, 14).Value = -18907.334
$ws.Cells.Item(58, 8).Value = 12827107
$ws.Cells.Item(58, 9).Value = 29413868
$ws.Cells.Item(58, 10).Value = 10065.046
$ws.Cells.Item(58, 11).Value = 29413868
$ws.Cells.Item(58, 12).Value = 10065.046
$ws.Cells.Item(58, 13).Value = -29413665
$ws.Cells.Item(58, 14).Value = -10471.046
$ws.Cells.Item(82, 8).Value = 24333.334
$ws.Cells.Item(82, 10).Value = 24333.334
$ws.Cells.Item(82, 12).Value = 24333.334
$ws.Cells.Item(82, 14).Value = -25055.334
$ws.Cells.Item(85, 8).Value = 24333.334
$ws.Cells.Item(85, 10).Value = 24333.334
$ws.Cells.Item(85, 12).Value = 24333.334
$ws.Cells.Item(85, 14).Value = -26829.334
$ws.Cells.Item(107, 8).Value = 1820.6666
$ws.Cells.Item(107, 9).Value = 668.3
$ws.Cells.Item(107, 10).Value = 2868.2727
$ws.Cells.Item(107, 11).Value = 668.3
$ws.Cells.Item(107, 12).Value = 2868.2727
$ws.Cells.Item(107, 13).Value = 1251.7
$ws.Cells.Item(107, 14).Value = -6708.2727
$ws.Cells.Item(132, 8).Value = 4968.8774
$ws.Cells.Item(132, 9).Value = 2105.4849
$ws.Cells.Item(132, 10).Value = 10874.625
$ws.Cells.Item(132, 11).Value = 6316.4547
$ws.Cells.Item(132, 12).Value = 32623.875
$ws.Cells.Item(132, 13).Value = -3786.4547
$ws.Cells.Item(132, 14).Value = -37683.875
$ws.Cells.Item(136, 8).Value = 12827107
$ws.Cells.Item(136, 9).Value = 29413868
$ws.Cells.Item(136, 10).Value = 10065.046
$ws.Cells.Item(136, 11).Value = 88241604
$ws.Cells.Item(136, 12).Value = 30195.138
$ws.Cells.Item(136, 13).Value = -88239054
$ws.Cells.Item(136, 14).Value = -35295.138

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 227.8077
$ws.Cells.Item(23, 9).Value = 160.8
$ws.Cells.Item(23, 10).Value = 269.6875
$ws.Cells.Item(23, 11).Value = 482.4
$ws.Cells.Item(23, 12).Value = 809.0625
$ws.Cells.Item(23, 13).Value = -247.4
$ws.Cells.Item(23, 14).Value = -1279.0625
$ws.Cells.Item(103, 8).Value = 1495.5
$ws.Cells.Item(103, 10).Value = 1893.7142
$ws.Cells.Item(103, 12).Value = 5681.142599999999
$ws.Cells.Item(103, 14).Value = -7439.142599999999
$ws.Cells.Item(128, 8).Value = 224484.5
$ws.Cells.Item(128, 9).Value = 224484.5
$ws.Cells.Item(128, 11).Value = 673453.5
$ws.Cells.Item(128, 13).Value = -668473.5
$ws.Cells.Item(131, 8).Value = 1908.5
$ws.Cells.Item(131, 9).Value = 1879.8334
$ws.Cells.Item(131, 10).Value = 1922.8334
$ws.Cells.Item(131, 11).Value = 5639.5002
$ws.Cells.Item(131, 12).Value = 5768.5002
$ws.Cells.Item(131, 13).Value = -599.5002000000004
$ws.Cells.Item(131, 14).Value = -15848.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 66.333336
$ws.Cells.Item(2, 9).Value = 52.75
$ws.Cells.Item(2, 10).Value = 175
$ws.Cells.Item(2, 11).Value = 52.75
$ws.Cells.Item(2, 12).Value = 175
$ws.Cells.Item(2, 13).Value = 60.25
$ws.Cells.Item(2, 14).Value = -401
$ws.Cells.Item(113, 8).Value = 291497.84
$ws.Cells.Item(113, 9).Value = 668969.8
$ws.Cells.Item(113, 10).Value = 8393.85
$ws.Cells.Item(113, 11).Value = 668969.8
$ws.Cells.Item(113, 12).Value = 8393.85
$ws.Cells.Item(113, 13).Value = -666799.8
$ws.Cells.Item(113, 14).Value = -12733.85
$ws.Cells.Item(122, 8).Value = 2423096.2
$ws.Cells.Item(122, 9).Value = 3460154.8
$ws.Cells.Item(122, 10).Value = 3293.2222
$ws.Cells.Item(122, 11).Value = 10380464.4
$ws.Cells.Item(122, 12).Value = 9879.6666
$ws.Cells.Item(122, 13).Value = -10378014.4
$ws.Cells.Item(122, 14).Value = -14779.6666
$ws.Cells.Item(140, 8).Value = 98780
$ws.Cells.Item(140, 10).Value = 98780
$ws.Cells.Item(140, 12).Value = 98780
$ws.Cells.Item(140, 14).Value = -109140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5350.1
$ws.Cells.Item(61, 9).Value = 1974.1428
$ws.Cells.Item(61, 11).Value = 1974.1428
$ws.Cells.Item(61, 13).Value = -1772.1428
$ws.Cells.Item(68, 8).Value = 5253.6924
$ws.Cells.Item(68, 9).Value = 4266.6665
$ws.Cells.Item(68, 11).Value = 4266.6665
$ws.Cells.Item(68, 13).Value = -3517.6665
$ws.Cells.Item(71, 8).Value = 5253.6924
$ws.Cells.Item(71, 9).Value = 4266.6665
$ws.Cells.Item(71, 11).Value = 21333.3325
$ws.Cells.Item(71, 13).Value = -17589.3325
$ws.Cells.Item(87, 8).Value = 56500
$ws.Cells.Item(87, 10).Value = 56500
$ws.Cells.Item(87, 12).Value = 56500
$ws.Cells.Item(87, 14).Value = -58746
$ws.Cells.Item(90, 8).Value = 56500
$ws.Cells.Item(90, 10).Value = 56500
$ws.Cells.Item(90, 12).Value = 169500
$ws.Cells.Item(90, 14).Value = -180732
$ws.Cells.Item(100, 8).Value = 3638.8696
$ws.Cells.Item(100, 10).Value = 5352.1113
$ws.Cells.Item(100, 12).Value = 5352.1113
$ws.Cells.Item(100, 14).Value = -6434.1113
$ws.Cells.Item(113, 8).Value = 5350.1
$ws.Cells.Item(113, 9).Value = 1974.1428
$ws.Cells.Item(113, 11).Value = 1974.1428
$ws.Cells.Item(113, 13).Value = 195.8571999999999
$ws.Cells.Item(122, 8).Value = 5486.081
$ws.Cells.Item(122, 9).Value = 4480
$ws.Cells.Item(122, 11).Value = 13440
$ws.Cells.Item(122, 13).Value = -10990
$ws.Cells.Item(136, 8).Value = 10599.917
$ws.Cells.Item(136, 9).Value = 2183
$ws.Cells.Item(136, 11).Value = 6549
$ws.Cells.Item(136, 13).Value = -3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 137012.97
$ws.Cells.Item(122, 9).Value = 212788.2
$ws.Cells.Item(122, 10).Value = 6128.4546
$ws.Cells.Item(122, 11).Value = 638364.6000000001
$ws.Cells.Item(122, 12).Value = 18385.3638
$ws.Cells.Item(122, 13).Value = -635914.6000000001
$ws.Cells.Item(122, 14).Value = -23285.3638
$ws.Cells.Item(126, 8).Value = 2192.2173
$ws.Cells.Item(126, 9).Value = 1834.4
$ws.Cells.Item(126, 11).Value = 5503.200000000001
$ws.Cells.Item(126, 13).Value = -3033.200000000001
$ws.Cells.Item(139, 8).Value = 83186.664
$ws.Cells.Item(139, 10).Value = 82335
$ws.Cells.Item(139, 12).Value = 82335
$ws.Cells.Item(139, 14).Value = -92615
